# RSTK-8093 - Inventory Location Adjustment
# Misc Change-CPQ CRUD operations
#
# The "Inventory Location details" sheet had a stray/duplicate CRUD test
# row (old row 11: multidiv serial / CON (Consigned) / KCB2T /
# Adjust Qty Down / serial down) removed; everything below it (old row 12)
# shifts up to take its place, which also drops the now-unused shared
# strings ("CON (Consigned)", "KCB2T", "serial down") from the workbook.

$wb = $excel.ActiveWorkbook

$wsItem = $wb.Worksheets.Item("Inventory Item")
$wsLoc  = $wb.Worksheets.Item("Inventory Location details")

# Remove the bad CRUD test row (row 11) from the location-details sheet.
$wsLoc.Rows("11:11").Delete() | Out-Null

# Leave the UI selection state the way the author's Excel session ended up:
# row 6 downward selected on "Inventory Item", and row 12 downward selected
# on "Inventory Location details" (just past the now-shorter data range).
$wsItem.Range("A6:A1048576").EntireRow.Select() | Out-Null
$wsLoc.Range("A12:A1048576").EntireRow.Select() | Out-Null
